$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append the new "control_ZMP1 / elpoli article / Apartado 5.1.4" row
# at the bottom of the table (§5.1.4 wording per the commit message).
$t.Rows.Add() | Out-Null
$idx = $t.Rows.Count
$t.Cell($idx, 1).Range.Text = "control_ZMP₁"
$t.Cell($idx, 2).Range.Text = "http://revistas.elpoli.edu.co/index.php/pol/article/viewFile/138/114"
$t.Cell($idx, 3).Range.Text = "Apartado 5.1.4"

# Table-level indent and default cell margin (tblInd -15 -> -20 dxa;
# tblCellMar left 93 -> 88 dxa). Word COM expresses these in points
# (1 pt = 20 dxa), so -20 dxa = -1.0 pt and 88 dxa = 4.4 pt.
$t.Rows.LeftIndent = -1.0
$t.LeftPadding = 4.4

# Every cell (the pre-existing ones plus the row just appended) also
# carries its own left-margin override that needs the same 93 -> 88 dxa
# (4.65 -> 4.4 pt) update.
foreach ($cell in $t.Range.Cells) {
    $cell.LeftPadding = 4.4
}
